$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '29.348.40'
$ws.Range('E2').Value = '  -0.24%  '
$ws.Range('D3').Value = '1.847.77'
$ws.Range('E3').Value = '  -0.22%  '
$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value = '0.9980'
$ws.Range('D4').ClearFormats()
$ws.Range('E4').Value = '  -0.19%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '240.09'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  -0.40%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '0.6269'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  -0.46%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.9989'
$ws.Range('D7').ClearFormats()
$ws.Range('E7').Value = '  -0.15%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.07597'
$ws.Range('D8').ClearFormats()
$ws.Range('E8').Value = '  -1.15%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.2906'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  -1.10%  '
$ws.Range('E10').Value = '  +0.26%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.07743'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  -0.13%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '5.020'
$ws.Range('D12').ClearFormats()
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '0.6790'
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  -0.34%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '0.00001050'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  -3.97%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '83.07'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  -0.76%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '6.123'
$ws.Range('D16').ClearFormats()
$ws.Range('E16').Value = '  -0.64%  '
$ws.Range('D17').Value = '29.392.55'
$ws.Range('E17').Value = '  -0.24%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '229.14'
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = '  -0.13%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '12.33'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  -1.22%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '0.9985'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  -0.21%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '7.461'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  +0.07%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '0.9973'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  -0.33%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '158.49'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  +0.91%  '
$ws.Range('E24').Value = '  -0.25%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '8.431'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  +0.49%  '
$ws.Range('E26').Value = '  -0.15%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '1.437'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  +9.28%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '1.466'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  -0.07%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '0.05597'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  -2.26%  '
$ws.Range('E30').Value = '  -0.73%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '4.067'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  +0.38%  '
$ws.Range('E32').Value = '  -0.33%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '1.830'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  -1.05%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '0.6959'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  -1.78%  '
$ws.Range('D36').Value = '1.233.49'
$ws.Range('E36').Value = '  +1.09%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '0.01798'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  +0.36%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '2.728'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  -1.78%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '6.377'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  -1.77%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '0.9042'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  -0.44%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '0.9988'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  -0.16%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '101.29'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  -0.37%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '65.42'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  -1.49%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '7.185'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  +0.63%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '0.3992'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  -0.65%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '9.013'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  +0.47%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '1.681'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  -0.15%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '0.1146'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  +1.16%  '
$ws.Range('E49').Value = '  -4.49%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '0.05696'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  -0.28%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '0.4621'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  -0.15%  '
